$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1557.1
$ws.Range("I58").Value = 187.5
$ws.Range("K58").Value = 562.5
$ws.Range("M58").Value = -412.5
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H132").Value = 2115.8125
$ws.Range("I132").Value = 1361.7142
$ws.Range("J132").Value = 7394.5
$ws.Range("K132").Value = 4085.1426
$ws.Range("L132").Value = 22183.5
$ws.Range("M132").Value = -1555.1426
$ws.Range("N132").Value = -27243.5
$ws.Range("H135").Value = 1179.4286
$ws.Range("J135").Value = 2000
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070
$ws.Range("H138").Value = 1883.3438
$ws.Range("I138").Value = 1194.2727
$ws.Range("K138").Value = 3582.8181
$ws.Range("M138").Value = 1557.1819

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4000.5
$ws.Range("I74").Value = 3501
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 3501
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -2627
$ws.Range("N74").Value = -6248
$ws.Range("H77").Value = 4000.5
$ws.Range("I77").Value = 3501
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 17505
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -13137
$ws.Range("N77").Value = -31236
$ws.Range("H122").Value = 2864
$ws.Range("I122").Value = 2262.2354
$ws.Range("K122").Value = 6786.706200000001
$ws.Range("M122").Value = -4336.706200000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4931.2856
$ws.Range("I20").Value = 3663.8
$ws.Range("K20").Value = 3663.8
$ws.Range("M20").Value = -3416.8
$ws.Range("H134").Value = 2430.4
$ws.Range("I134").Value = 2413
$ws.Range("K134").Value = 7239
$ws.Range("M134").Value = -4704

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2777.5833
$ws.Range("I122").Value = 2208.875
$ws.Range("K122").Value = 6626.625
$ws.Range("M122").Value = -4176.625
$ws.Range("H132").Value = 6364.048
$ws.Range("J132").Value = 6215.231
$ws.Range("L132").Value = 18645.693
$ws.Range("N132").Value = -23705.693
$ws.Range("H134").Value = 2722.3635
$ws.Range("I134").Value = 2033.1666
$ws.Range("J134").Value = 3549.4
$ws.Range("K134").Value = 6099.4998
$ws.Range("L134").Value = 10648.2
$ws.Range("M134").Value = -3564.4998
$ws.Range("N134").Value = -15718.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1861.4634
$ws.Range("I11").Value = 1568.3334
$ws.Range("J11").Value = 1982.7587
$ws.Range("K11").Value = 4705.0002
$ws.Range("L11").Value = 5948.2761
$ws.Range("M11").Value = -4565.0002
$ws.Range("N11").Value = -6228.2761
$ws.Range("H46").Value = 5000400
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H59").Value = 2502.5
$ws.Range("I59").Value = 1005
$ws.Range("J59").Value = 4000
$ws.Range("K59").Value = 3015
$ws.Range("L59").Value = 12000
$ws.Range("M59").Value = -2475
$ws.Range("N59").Value = -13080
$ws.Range("H99").Value = 1990
$ws.Range("I99").Value = 460.4
$ws.Range("J99").Value = 3902
$ws.Range("K99").Value = 1381.2
$ws.Range("L99").Value = 11706
$ws.Range("M99").Value = 864.8000000000002
$ws.Range("N99").Value = -16198
$ws.Range("H102").Value = 8996.75
$ws.Range("J102").Value = 8996.75
$ws.Range("L102").Value = 26990.25
$ws.Range("N102").Value = -31858.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H70").Value = 10576.846
$ws.Range("I70").Value = 7749.4
$ws.Range("K70").Value = 7749.4
$ws.Range("M70").Value = -7479.4
$ws.Range("H73").Value = 10576.846
$ws.Range("I73").Value = 7749.4
$ws.Range("K73").Value = 7749.4
$ws.Range("M73").Value = -6813.4
$ws.Range("H93").Value = 28998.5
$ws.Range("J93").Value = 28998.5
$ws.Range("L93").Value = 28998.5
$ws.Range("N93").Value = -32742.5
$ws.Range("H122").Value = 37046.645
$ws.Range("I122").Value = 3683.1904
$ws.Range("J122").Value = 107109.9
$ws.Range("K122").Value = 11049.5712
$ws.Range("L122").Value = 321329.7
$ws.Range("M122").Value = -8599.5712
$ws.Range("N122").Value = -326229.7
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -79900

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2003
$ws.Range("I93").Value = 2003
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2003
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -755
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 3899.4849
$ws.Range("I122").Value = 2789.9
$ws.Range("J122").Value = 4381.913
$ws.Range("K122").Value = 8369.700000000001
$ws.Range("L122").Value = 13145.739
$ws.Range("M122").Value = -5919.700000000001
$ws.Range("N122").Value = -18045.739
$ws.Range("H136").Value = 5050
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 34000
$ws.Range("I40").Value = 34000
$ws.Range("K40").Value = 34000
$ws.Range("M40").Value = -33851
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H81").Value = 1762.25
$ws.Range("I81").Value = 1856.8572
$ws.Range("K81").Value = 3713.7144
$ws.Range("M81").Value = -2652.7144
$ws.Range("H84").Value = 1762.25
$ws.Range("I84").Value = 1856.8572
$ws.Range("K84").Value = 18568.572
$ws.Range("M84").Value = -13264.572
$ws.Range("H122").Value = 1667.5883
$ws.Range("I122").Value = 1546.6
$ws.Range("K122").Value = 4639.799999999999
$ws.Range("M122").Value = -2189.799999999999
$ws.Range("H132").Value = 2416
$ws.Range("I132").Value = 2392.1428
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 7176.428400000001
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -4646.428400000001
$ws.Range("N132").Value = -12558.5
$ws.Range("H135").Value = 37905
$ws.Range("J135").Value = 37905
$ws.Range("L135").Value = 37905
$ws.Range("N135").Value = -48045
$ws.Range("H136").Value = 6789.4
$ws.Range("I136").Value = 6789.4
$ws.Range("K136").Value = 20368.2
$ws.Range("M136").Value = -17818.2
